# Apply the dated worksheet update: change the date and all division
# problems in the document, per the commit "Update master to output
# generated at c8c62b6".

$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-06-22 Sunday"; new = "2025-06-23 Monday"},
    @{old = "561÷8="; new = "818÷2="},
    @{old = "676÷4="; new = "288÷6="},
    @{old = "751÷6="; new = "828÷5="},
    @{old = "933÷9="; new = "541÷8="},
    @{old = "152÷2="; new = "230÷4="},
    @{old = "753÷2="; new = "557÷6="},
    @{old = "710÷6="; new = "884÷7="},
    @{old = "958÷7="; new = "375÷3="},
    @{old = "781÷3="; new = "891÷7="},
    @{old = "540÷2="; new = "137÷7="},
    @{old = "495÷7="; new = "793÷9="},
    @{old = "836÷9="; new = "714÷5="},
    @{old = "557÷8="; new = "674÷9="},
    @{old = "285÷4="; new = "161÷2="},
    @{old = "266÷9="; new = "704÷2="},
    @{old = "633÷5="; new = "682÷9="},
    @{old = "182÷9="; new = "185÷4="},
    @{old = "584÷9="; new = "797÷4="},
    @{old = "461÷8="; new = "318÷8="},
    @{old = "485÷6="; new = "706÷9="},
    @{old = "518÷4="; new = "702÷5="},
    @{old = "553÷7="; new = "351÷6="},
    @{old = "627÷3="; new = "784÷4="},
    @{old = "371÷2="; new = "469÷9="},
    @{old = "472÷9="; new = "562÷2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
